$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.357.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +1.16%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.667.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +1.02%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +1.00%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''219.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +1.15%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.5352'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +1.19%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.92%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.2665'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +2.60%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.06394'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +1.33%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''20.90'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +2.97%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.07841'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.63%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''4.562'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +1.03%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.660.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +0.09%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''1.896.43'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +0.98%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.5546'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.41%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''0.0₅8197'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +0.19%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  +1.01%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''26.384.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +1.25%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.92%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''4.685'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +2.49%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''195.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +2.56%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''10.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +2.02%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''6.043'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +0.55%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''  +0.95%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''146.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +2.24%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.30%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''7.234'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.40%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''16.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +0.18%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''1.504'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +4.18%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''0.05855'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.82%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''1.285'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +1.03%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''3.586'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +1.29%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''3.293'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +1.15%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  +1.41%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''0.9727'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +3.25%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  +1.47%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''2.423'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = '''0.5831'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +1.53%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.01609'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +0.26%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''1.076.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +4.58%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.8646'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +1.80%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''5.861'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +2.57%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  +0.92%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''104.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -0.20%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''1.806.16'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.74%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''58.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +1.95%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''1.015'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +1.38%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''0.4391'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +1.46%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''8.026'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +2.28%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.0₈102'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -8.49%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '''  +0.61%  '
$ws.Range('E51').Style = 'Normal'
